$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2248.7556
$ws.Range("I15").Value = 2248.7556
$ws.Range("K15").Value = 6746.266799999999
$ws.Range("M15").Value = -6577.266799999999

$ws.Range("H19").Value = 995.94116
$ws.Range("J19").Value = 873.7
$ws.Range("L19").Value = 873.7
$ws.Range("N19").Value = -1223.7

$ws.Range("H58").Value = 1300.2727
$ws.Range("I58").Value = 701
$ws.Range("K58").Value = 2103
$ws.Range("M58").Value = -1953

$ws.Range("H106").Value = 2186.52
$ws.Range("I106").Value = 2186.52
$ws.Range("K106").Value = 2186.52
$ws.Range("M106").Value = -1555.52

$ws.Range("H135").Value = 41514.88
$ws.Range("I135").Value = 1253.2
$ws.Range("J135").Value = 202561.6
$ws.Range("K135").Value = 11278.8
$ws.Range("L135").Value = 1823054.4
$ws.Range("M135").Value = -8743.800000000001
$ws.Range("N135").Value = -1828124.4

$ws.Range("H138").Value = 3935.4
$ws.Range("J138").Value = 4686.0415
$ws.Range("L138").Value = 14058.1245
$ws.Range("N138").Value = -24338.1245

$ws.Range("H141").Value = 9549.852999999999
$ws.Range("I141").Value = 9026.6
$ws.Range("J141").Value = 13474.25
$ws.Range("K141").Value = 27079.8
$ws.Range("L141").Value = 40422.75
$ws.Range("M141").Value = -21899.8
$ws.Range("N141").Value = -50782.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 429.75
$ws.Range("I5").Value = 457.91666
$ws.Range("K5").Value = 457.91666
$ws.Range("M5").Value = -345.91666

$ws.Range("H26").Value = 4336
$ws.Range("I26").Value = 4336
$ws.Range("K26").Value = 4336
$ws.Range("M26").Value = -4006

$ws.Range("H32").Value = 16887.096
$ws.Range("I32").Value = 14645.156
$ws.Range("J32").Value = 20474.2
$ws.Range("K32").Value = 14645.156
$ws.Range("L32").Value = 20474.2
$ws.Range("M32").Value = -14358.156
$ws.Range("N32").Value = -21048.2

$ws.Range("H45").Value = 420976.75
$ws.Range("I45").Value = 1252671.2
$ws.Range("K45").Value = 1252671.2
$ws.Range("M45").Value = -1252294.2

$ws.Range("H61").Value = 4145.1665
$ws.Range("I61").Value = 3919.5469
$ws.Range("J61").Value = 5950.125
$ws.Range("K61").Value = 3919.5469
$ws.Range("L61").Value = 5950.125
$ws.Range("M61").Value = -3707.5469
$ws.Range("N61").Value = -6374.125

$ws.Range("H110").Value = 8337
$ws.Range("I110").Value = 7505.5
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 7505.5
$ws.Range("L110").Value = 10000
$ws.Range("M110").Value = -5460.5
$ws.Range("N110").Value = -14090

$ws.Range("H132").Value = 18752.781
$ws.Range("I132").Value = 21456.111
$ws.Range("K132").Value = 64368.333
$ws.Range("M132").Value = -61838.333

$ws.Range("H136").Value = 4145.1665
$ws.Range("I136").Value = 3919.5469
$ws.Range("J136").Value = 5950.125
$ws.Range("K136").Value = 11758.6407
$ws.Range("L136").Value = 17850.375
$ws.Range("M136").Value = -9208.6407
$ws.Range("N136").Value = -22950.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 429.75
$ws.Range("I4").Value = 457.91666
$ws.Range("K4").Value = 457.91666
$ws.Range("M4").Value = -342.91666

$ws.Range("H21").Value = 27150
$ws.Range("J21").Value = 27150
$ws.Range("L21").Value = 27150
$ws.Range("N21").Value = -27622

$ws.Range("H59").Value = 121262
$ws.Range("J59").Value = 121262
$ws.Range("L59").Value = 121262
$ws.Range("N59").Value = -122956

$ws.Range("H125").Value = 95313.8
$ws.Range("J125").Value = 95313.8
$ws.Range("L125").Value = 95313.8
$ws.Range("N125").Value = -105153.8

$ws.Range("H134").Value = 2363.5
$ws.Range("I134").Value = 2220.653
$ws.Range("J134").Value = 4696.6665
$ws.Range("K134").Value = 6661.958999999999
$ws.Range("L134").Value = 14089.9995
$ws.Range("M134").Value = -4126.958999999999
$ws.Range("N134").Value = -19159.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()

$ws.Range("H31").Value = 7612.2
$ws.Range("I31").Value = 4766.6665
$ws.Range("J31").Value = 8831.714
$ws.Range("K31").Value = 4766.6665
$ws.Range("L31").Value = 8831.714
$ws.Range("M31").Value = -4471.6665
$ws.Range("N31").Value = -9421.714

$ws.Range("H34").Value = 7612.2
$ws.Range("I34").Value = 4766.6665
$ws.Range("J34").Value = 8831.714
$ws.Range("K34").Value = 4766.6665
$ws.Range("L34").Value = 8831.714
$ws.Range("M34").Value = -4564.6665
$ws.Range("N34").Value = -9235.714

$ws.Range("H58").Value = 94592.63
$ws.Range("I58").Value = 103810.5
$ws.Range("J58").Value = 2414
$ws.Range("K58").Value = 103810.5
$ws.Range("L58").Value = 2414
$ws.Range("M58").Value = -103607.5
$ws.Range("N58").Value = -2820

$ws.Range("H99").Value = 5070
$ws.Range("I99").Value = 3443.75
$ws.Range("K99").Value = 3443.75
$ws.Range("M99").Value = -1945.75

$ws.Range("H126").Value = 5070
$ws.Range("I126").Value = 3443.75
$ws.Range("K126").Value = 10331.25
$ws.Range("M126").Value = -7861.25

$ws.Range("H136").Value = 94592.63
$ws.Range("I136").Value = 103810.5
$ws.Range("J136").Value = 2414
$ws.Range("K136").Value = 311431.5
$ws.Range("L136").Value = 7242
$ws.Range("M136").Value = -308881.5
$ws.Range("N136").Value = -12342

$ws.Range("H140").Value = 79995.664
$ws.Range("J140").Value = 79995.664
$ws.Range("L140").Value = 79995.664
$ws.Range("N140").Value = -90355.664

$ws.Range("H141").Value = 401429.8
$ws.Range("J141").Value = 401429.8
$ws.Range("L141").Value = 401429.8
$ws.Range("N141").Value = -411789.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2980.842
$ws.Range("J5").Value = 3778.75
$ws.Range("L5").Value = 11336.25
$ws.Range("N5").Value = -11560.25

$ws.Range("H107").Value = 557.53845
$ws.Range("I107").Value = 750
$ws.Range("J107").Value = 541.5
$ws.Range("K107").Value = 2250
$ws.Range("L107").Value = 1624.5
$ws.Range("M107").Value = -330
$ws.Range("N107").Value = -5464.5

$ws.Range("H113").Value = 2770.5715
$ws.Range("J113").Value = 2770.5715
$ws.Range("L113").Value = 8311.7145
$ws.Range("N113").Value = -12651.7145

$ws.Range("H135").Value = 2980.842
$ws.Range("J135").Value = 3778.75
$ws.Range("L135").Value = 34008.75
$ws.Range("N135").Value = -39078.75

$ws.Range("H137").Value = 3390.5833
$ws.Range("J137").Value = 3889.5715
$ws.Range("L137").Value = 11668.7145
$ws.Range("N137").Value = -21868.7145

$ws.Range("H138").Value = 3052.8572
$ws.Range("I138").Value = 3145
$ws.Range("J138").Value = 2500
$ws.Range("K138").Value = 9435
$ws.Range("L138").Value = 7500
$ws.Range("M138").Value = -4295
$ws.Range("N138").Value = -17780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1225
$ws.Range("J3").Value = 1666.6666
$ws.Range("L3").Value = 1666.6666
$ws.Range("N3").Value = -1898.6666

$ws.Range("H11").Value = 300
$ws.Range("I11").Value = 300
$ws.Range("K11").Value = 300
$ws.Range("M11").Value = -161

$ws.Range("H43").Value = 17901.25
$ws.Range("I43").Value = 3963.2
$ws.Range("J43").Value = 27857
$ws.Range("K43").Value = 3963.2
$ws.Range("L43").Value = 27857
$ws.Range("M43").Value = -3812.2
$ws.Range("N43").Value = -28159

$ws.Range("H122").Value = 4276.2583
$ws.Range("I122").Value = 3106.9
$ws.Range("J122").Value = 4833.095
$ws.Range("K122").Value = 9320.700000000001
$ws.Range("L122").Value = 14499.285
$ws.Range("M122").Value = -6870.700000000001
$ws.Range("N122").Value = -19399.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 366221.56
$ws.Range("I7").Value = 409484.16
$ws.Range("J7").Value = 5700
$ws.Range("K7").Value = 409484.16
$ws.Range("L7").Value = 5700
$ws.Range("M7").Value = -409372.16
$ws.Range("N7").Value = -5924

$ws.Range("H22").Value = 2512.8667
$ws.Range("I22").Value = 1426
$ws.Range("J22").Value = 3755
$ws.Range("K22").Value = 1426
$ws.Range("L22").Value = 3755
$ws.Range("M22").Value = -1131
$ws.Range("N22").Value = -4345

$ws.Range("H27").Value = 2512.8667
$ws.Range("I27").Value = 1426
$ws.Range("J27").Value = 3755
$ws.Range("K27").Value = 1426
$ws.Range("L27").Value = 3755
$ws.Range("M27").Value = -1319
$ws.Range("N27").Value = -3969

$ws.Range("H61").Value = 2655
$ws.Range("I61").Value = 2554.8823
$ws.Range("J61").Value = 3222.3333
$ws.Range("K61").Value = 2554.8823
$ws.Range("L61").Value = 3222.3333
$ws.Range("M61").Value = -2352.8823
$ws.Range("N61").Value = -3626.3333

$ws.Range("H100").Value = 2920.3333
$ws.Range("I100").Value = 2333.6
$ws.Range("K100").Value = 2333.6
$ws.Range("M100").Value = -1792.6

$ws.Range("H113").Value = 2655
$ws.Range("I113").Value = 2554.8823
$ws.Range("J113").Value = 3222.3333
$ws.Range("K113").Value = 2554.8823
$ws.Range("L113").Value = 3222.3333
$ws.Range("M113").Value = -384.8823000000002
$ws.Range("N113").Value = -7562.3333

$ws.Range("H122").Value = 4724.1934
$ws.Range("I122").Value = 4389.2856
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 13167.8568
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -10717.8568
$ws.Range("N122").Value = -19900

$ws.Range("H126").Value = 366221.56
$ws.Range("I126").Value = 409484.16
$ws.Range("J126").Value = 5700
$ws.Range("K126").Value = 1228452.48
$ws.Range("L126").Value = 17100
$ws.Range("M126").Value = -1225982.48
$ws.Range("N126").Value = -22040

$ws.Range("H136").Value = 7473.5
$ws.Range("I136").Value = 2899
$ws.Range("K136").Value = 8697
$ws.Range("M136").Value = -6147

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3170.7026
$ws.Range("I136").Value = 2782.303
$ws.Range("J136").Value = 6375
$ws.Range("K136").Value = 8346.909
$ws.Range("L136").Value = 19125
$ws.Range("M136").Value = -5796.909
$ws.Range("N136").Value = -24225
